$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value2 = 3812.375
$ws.Range("I74").Value2 = 3699.7778
$ws.Range("J74").Value2 = 3957.1428
$ws.Range("K74").Value2 = 3699.7778
$ws.Range("L74").Value2 = 3957.1428
$ws.Range("M74").Value2 = -2763.7778
$ws.Range("N74").Value2 = -5829.1428
$ws.Range("H77").Value2 = 3812.375
$ws.Range("I77").Value2 = 3699.7778
$ws.Range("J77").Value2 = 3957.1428
$ws.Range("K77").Value2 = 18498.889
$ws.Range("L77").Value2 = 19785.714
$ws.Range("M77").Value2 = -13818.889
$ws.Range("N77").Value2 = -29145.714
$ws.Range("H132").Value2 = 344204.22
$ws.Range("I132").Value2 = 434891.06
$ws.Range("K132").Value2 = 1304673.18
$ws.Range("M132").Value2 = -1302143.18
$ws.Range("H136").Value2 = 48500
$ws.Range("J136").Value2 = 48500
$ws.Range("L136").Value2 = 48500
$ws.Range("N136").Value2 = -58700
$ws.Range("H137").Value2 = 33334946
$ws.Range("I137").Value2 = 47620156
$ws.Range("K137").Value2 = 142860468
$ws.Range("M137").Value2 = -142857918
$ws.Range("H139").Value2 = 42333.332
$ws.Range("J139").Value2 = 42333.332
$ws.Range("L139").Value2 = 42333.332
$ws.Range("N139").Value2 = -52613.332
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value2 = 43799.4
$ws.Range("J80").Value2 = 29999.25
$ws.Range("L80").Value2 = 29999.25
$ws.Range("N80").Value2 = -31995.25
$ws.Range("H83").Value2 = 43799.4
$ws.Range("J83").Value2 = 29999.25
$ws.Range("L83").Value2 = 89997.75
$ws.Range("N83").Value2 = -99981.75
$ws.Range("H132").Value2 = 2979.8293
$ws.Range("I132").Value2 = 2710.8147
$ws.Range("K132").Value2 = 8132.4441
$ws.Range("M132").Value2 = -5602.4441
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value2 = 873.25
$ws.Range("I64").Value2 = 876.5
$ws.Range("J64").Value2 = 870
$ws.Range("K64").Value2 = 876.5
$ws.Range("L64").Value2 = 870
$ws.Range("M64").Value2 = -651.5
$ws.Range("N64").Value2 = -1320
$ws.Range("H67").Value2 = 873.25
$ws.Range("I67").Value2 = 876.5
$ws.Range("J67").Value2 = 870
$ws.Range("K67").Value2 = 876.5
$ws.Range("L67").Value2 = 870
$ws.Range("M67").Value2 = -96.5
$ws.Range("N67").Value2 = -2430
$ws.Range("H86").Value2 = 22141.4
$ws.Range("I86").Value2 = 2675
$ws.Range("K86").Value2 = 2675
$ws.Range("M86").Value2 = -1552
$ws.Range("H89").Value2 = 22141.4
$ws.Range("I89").Value2 = 2675
$ws.Range("K89").Value2 = 13375
$ws.Range("M89").Value2 = -7759
$ws.Range("H135").Value2 = 56920
$ws.Range("J135").Value2 = 56920
$ws.Range("L135").Value2 = 56920
$ws.Range("N135").Value2 = -67060
$ws.Range("H138").Value2 = 44500
$ws.Range("J138").Value2 = 44500
$ws.Range("L138").Value2 = 44500
$ws.Range("N138").Value2 = -54780
$ws.Range("H140").Value2 = 37926.668
$ws.Range("J140").Value2 = 37926.668
$ws.Range("L140").Value2 = 37926.668
$ws.Range("N140").Value2 = -48286.668
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1400.75
$ws.Range("I16").Value2 = 1334.3334
$ws.Range("J16").Value2 = 1600
$ws.Range("K16").Value2 = 1334.3334
$ws.Range("L16").Value2 = 1600
$ws.Range("M16").Value2 = -1047.3334
$ws.Range("N16").Value2 = -2174
$ws.Range("H22").Value2 = 794.2308
$ws.Range("I22").Value2 = 222.875
$ws.Range("J22").Value2 = 1708.4
$ws.Range("K22").Value2 = 222.875
$ws.Range("L22").Value2 = 1708.4
$ws.Range("M22").Value2 = 127.125
$ws.Range("N22").Value2 = -2408.4
$ws.Range("H31").Value2 = 1104.0312
$ws.Range("I31").Value2 = 960.9666999999999
$ws.Range("J31").Value2 = 3250
$ws.Range("K31").Value2 = 960.9666999999999
$ws.Range("L31").Value2 = 3250
$ws.Range("M31").Value2 = -665.9666999999999
$ws.Range("N31").Value2 = -3840
$ws.Range("H34").Value2 = 1104.0312
$ws.Range("I34").Value2 = 960.9666999999999
$ws.Range("J34").Value2 = 3250
$ws.Range("K34").Value2 = 960.9666999999999
$ws.Range("L34").Value2 = 3250
$ws.Range("M34").Value2 = -758.9666999999999
$ws.Range("N34").Value2 = -3654
$ws.Range("H62").Value2 = 19046.867
$ws.Range("I62").Value2 = 51975
$ws.Range("K62").Value2 = 51975
$ws.Range("M62").Value2 = -51351
$ws.Range("H65").Value2 = 19046.867
$ws.Range("I65").Value2 = 51975
$ws.Range("K65").Value2 = 259875
$ws.Range("M65").Value2 = -256755
$ws.Range("H99").Value2 = 6251297.5
$ws.Range("I99").Value2 = 8929632
$ws.Range("K99").Value2 = 8929632
$ws.Range("M99").Value2 = -8928134
$ws.Range("H105").Value2 = 681.2222
$ws.Range("I105").Value2 = 616.375
$ws.Range("J105").Value2 = 1200
$ws.Range("K105").Value2 = 616.375
$ws.Range("L105").Value2 = 1200
$ws.Range("M105").Value2 = 1130.625
$ws.Range("N105").Value2 = -4694
$ws.Range("H113").Value2 = 1400.75
$ws.Range("I113").Value2 = 1334.3334
$ws.Range("J113").Value2 = 1600
$ws.Range("K113").Value2 = 1334.3334
$ws.Range("L113").Value2 = 1600
$ws.Range("M113").Value2 = 835.6666
$ws.Range("N113").Value2 = -5940
$ws.Range("H115").Value2 = 27929.5
$ws.Range("J115").Value2 = 27929.5
$ws.Range("L115").Value2 = 27929.5
$ws.Range("N115").Value2 = -30279.5
$ws.Range("H126").Value2 = 6251297.5
$ws.Range("I126").Value2 = 8929632
$ws.Range("K126").Value2 = 26788896
$ws.Range("M126").Value2 = -26786426
$ws.Range("H134").Value2 = 2320
$ws.Range("I134").Value2 = 1184.9656
$ws.Range("J134").Value2 = 5312.364
$ws.Range("K134").Value2 = 3554.8968
$ws.Range("L134").Value2 = 15937.092
$ws.Range("M134").Value2 = -1019.8968
$ws.Range("N134").Value2 = -21007.092
$ws.Range("H137").Value2 = 34260
$ws.Range("J137").Value2 = 46390
$ws.Range("L137").Value2 = 46390
$ws.Range("N137").Value2 = -56590
$ws.Range("H138").Value2 = 44500
$ws.Range("J138").Value2 = 44500
$ws.Range("L138").Value2 = 44500
$ws.Range("N138").Value2 = -54780
$ws.Range("H140").Value2 = 45780
$ws.Range("J140").Value2 = 45780
$ws.Range("L140").Value2 = 45780
$ws.Range("N140").Value2 = -56140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value2 = 10204819
$ws.Range("I113").Value2 = 576
$ws.Range("J113").Value2 = 13158679
$ws.Range("K113").Value2 = 1728
$ws.Range("L113").Value2 = 39476037
$ws.Range("M113").Value2 = 442
$ws.Range("N113").Value2 = -39480377
$ws.Range("H132").Value2 = 1083.0385
$ws.Range("J132").Value2 = 1312.2
$ws.Range("L132").Value2 = 11809.8
$ws.Range("N132").Value2 = -16869.8
$ws.Range("H134").Value2 = 5210.909
$ws.Range("I134").Value2 = 3169.5715
$ws.Range("J134").Value2 = 8783.25
$ws.Range("K134").Value2 = 9508.7145
$ws.Range("L134").Value2 = 26349.75
$ws.Range("M134").Value2 = -4438.7145
$ws.Range("N134").Value2 = -36489.75
$ws.Range("H136").Value2 = 2718.4546
$ws.Range("I136").Value2 = 2157.25
$ws.Range("J136").Value2 = 2795.862
$ws.Range("K136").Value2 = 6471.75
$ws.Range("L136").Value2 = 8387.585999999999
$ws.Range("M136").Value2 = -1371.75
$ws.Range("N136").Value2 = -18587.586
$ws.Range("H139").Value2 = 2059.0356
$ws.Range("I139").Value2 = 1639.7084
$ws.Range("K139").Value2 = 4919.1252
$ws.Range("M139").Value2 = 220.8747999999996
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 1327.25
$ws.Range("I97").Value2 = 1327.25
$ws.Range("J97").Value2 = 0
$ws.Range("K97").Value2 = 1327.25
$ws.Range("L97").Value2 = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value2 = -831.25
$ws.Range("H113").Value2 = 2000
$ws.Range("I113").Value2 = 2000
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 2000
$ws.Range("L113").Value2 = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value2 = 170
$ws.Range("H140").Value2 = 50640
$ws.Range("J140").Value2 = 50640
$ws.Range("L140").Value2 = 50640
$ws.Range("N140").Value2 = -61000
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 1126.8
$ws.Range("I46").Value2 = 1040
$ws.Range("J46").Value2 = 1300.4
$ws.Range("K46").Value2 = 1040
$ws.Range("L46").Value2 = 1300.4
$ws.Range("M46").Value2 = -852
$ws.Range("N46").Value2 = -1676.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value2 = 33695
$ws.Range("J57").Value2 = 33695
$ws.Range("L57").Value2 = 33695
$ws.Range("N57").Value2 = -35203
$ws.Range("H74").Value2 = 17853.334
$ws.Range("I74").Value2 = 12569
$ws.Range("J74").Value2 = 20495.5
$ws.Range("K74").Value2 = 12569
$ws.Range("L74").Value2 = 20495.5
$ws.Range("M74").Value2 = -11633
$ws.Range("N74").Value2 = -22367.5
$ws.Range("H77").Value2 = 17853.334
$ws.Range("I77").Value2 = 12569
$ws.Range("J77").Value2 = 20495.5
$ws.Range("K77").Value2 = 37707
$ws.Range("L77").Value2 = 61486.5
$ws.Range("M77").Value2 = -33027
$ws.Range("N77").Value2 = -70846.5
$ws.Range("H100").Value2 = 462.75
$ws.Range("I100").Value2 = 350.33334
$ws.Range("J100").Value2 = 800
$ws.Range("K100").Value2 = 700.66668
$ws.Range("L100").Value2 = 1600
$ws.Range("M100").Value2 = -159.66668
$ws.Range("N100").Value2 = -2682
$ws.Range("H136").Value2 = 14538147
$ws.Range("I136").Value2 = 20897368
$ws.Range("J136").Value2 = 2787.4285
$ws.Range("K136").Value2 = 62692104
$ws.Range("L136").Value2 = 8362.2855
$ws.Range("M136").Value2 = -62689554
$ws.Range("N136").Value2 = -13462.2855
$ws.Range("H137").Value2 = 43750
$ws.Range("J137").Value2 = 43750
$ws.Range("L137").Value2 = 43750
$ws.Range("N137").Value2 = -53950
$ws.Range("H139").Value2 = 53500
$ws.Range("J139").Value2 = 53500
$ws.Range("L139").Value2 = 53500
$ws.Range("N139").Value2 = -63780
